$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.938.97"
$ws.Range("E2").Value = "  +0.58%  "
$ws.Range("D3").Value = "1.811.08"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("E5").Value = "  -0.13%  "
$ws.Range("D6").Value = "1.001"
$ws.Range("D7").Value = "0.4989"
$ws.Range("E7").Value = "  -2.47%  "
$ws.Range("D8").Value = "0.3904"
$ws.Range("E8").Value = "  +1.17%  "
$ws.Range("D9").Value = "0.09782"
$ws.Range("E9").Value = "  +25.15%  "
$ws.Range("E10").Value = "  +0.86%  "
$ws.Range("D11").Value = "40.83"
$ws.Range("E11").Value = "  +0.36%  "
$ws.Range("D12").Value = "6.409"
$ws.Range("D13").Value = "1.001"
$ws.Range("E13").Value = "  -0.08%  "
$ws.Range("E14").Value = "  +1.26%  "
$ws.Range("D15").Value = "1.811.35"
$ws.Range("E15").Value = "  +2.09%  "
$ws.Range("D16").Value = "7.258"
$ws.Range("E16").Value = "  +0.62%  "
$ws.Range("D18").Value = "92.30"
$ws.Range("E18").Value = "  +0.87%  "
$ws.Range("D19").Value = "0.06637"
$ws.Range("E19").Value = "  +1.23%  "
$ws.Range("D20").Value = "1.000"
$ws.Range("E20").Value = "  -0.13%  "
$ws.Range("D21").Value = "17.14"
$ws.Range("E21").Value = "  +0.63%  "
$ws.Range("D22").Value = "5.897"
$ws.Range("E22").Value = "  -0.18%  "
$ws.Range("D23").Value = "27.998.06"
$ws.Range("E23").Value = "  +0.61%  "
$ws.Range("D24").Value = "11.07"
$ws.Range("E24").Value = "  +0.64%  "
$ws.Range("D25").Value = "2.244"
$ws.Range("E25").Value = "  +0.75%  "
$ws.Range("D26").Value = "158.43"
$ws.Range("E26").Value = "  -0.92%  "
$ws.Range("D27").Value = "2.020.63"
$ws.Range("E27").Value = "  +1.84%  "
$ws.Range("D29").Value = "2.382"
$ws.Range("E29").Value = "  +0.53%  "
$ws.Range("D30").Value = "126.55"
$ws.Range("E30").Value = "  +2.63%  "
$ws.Range("E31").Value = "  -1.07%  "
$ws.Range("E32").Value = "  -0.36%  "
$ws.Range("D33").Value = "5.547"
$ws.Range("E33").Value = "  +1.16%  "
$ws.Range("D34").Value = "3.595"
$ws.Range("E34").Value = "  -1.10%  "
$ws.Range("D35").Value = "0.06713"
$ws.Range("E35").Value = "  -5.08%  "
$ws.Range("E36").Value = "  +0.84%  "
$ws.Range("D37").Value = "8.837"
$ws.Range("E37").Value = "  +0.38%  "
$ws.Range("D38").Value = "0.2135"
$ws.Range("D39").Value = "4.916"
$ws.Range("E39").Value = "  -1.64%  "
$ws.Range("E40").Value = "  -1.91%  "
$ws.Range("D41").Value = "0.6155"
$ws.Range("E41").Value = "  +1.14%  "
$ws.Range("E42").Value = "  +1.65%  "
$ws.Range("D43").Value = "1.000"
$ws.Range("E43").Value = "  -0.07%  "
$ws.Range("D44").Value = "13.17"
$ws.Range("E44").Value = "  -0.07%  "
$ws.Range("D45").Value = "0.5877"
$ws.Range("E45").Value = "  -0.30%  "
$ws.Range("D46").Value = "1.284"
$ws.Range("E46").Value = "  -2.69%  "
$ws.Range("E47").Value = "  -0.39%  "
$ws.Range("D48").Value = "123.34"
$ws.Range("E48").Value = "  -2.00%  "
$ws.Range("D49").Value = "1.932"
$ws.Range("E49").Value = "  +1.86%  "
$ws.Range("D50").Value = "1.177"
$ws.Range("E50").Value = "  -2.02%  "
$ws.Range("D51").Value = "0.06766"
$ws.Range("E51").Value = "  -1.24%  "
